# Updates the crypto price/volume table on Sheet1 with freshly scraped
# values (per the "Updated cryptos list ... with GitHub Actions" commit).
# Numeric-looking Price (column D) values are written with a leading
# apostrophe so Excel stores them as text (matching the original
# inlineStr cells) instead of auto-converting them to numbers and
# dropping significant trailing zeros (e.g. "19.20" -> 19.2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.800.56'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.637.89'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = '''218.67'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').Value = '''0.0621'
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').Value = '''19.20'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = '''0.0844'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').Value = '1.867.11'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '1.632.89'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = '''4.13'
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '''0.523'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '''64.65'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').Value = '26.794.27'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').Value = '''214.44'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = '''4.35'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '''6.31'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('D24').Value = '''9.09'
$ws.Range('E24').Value = '  -2.79%  '
$ws.Range('D25').Value = '''147.38'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').Value = '''7.02'
$ws.Range('E28').Value = '  -1.20%  '
$ws.Range('D29').Value = '''15.65'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '''0.0504'
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').Value = '''2.97'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = '1.259.44'
$ws.Range('E35').Value = '  -2.07%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('D37').Value = '''0.0174'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').Value = '''0.526'
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('D39').Value = '''0.812'
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('D41').Value = '''0.805'
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('D43').Value = '1.778.44'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('D44').Value = '''2.13'
$ws.Range('E44').Value = '  -4.53%  '
$ws.Range('D45').Value = '''91.98'
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('D46').Value = '''59.98'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.0962'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').Value = '''1.01'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.405'
$ws.Range('E51').Value = '  -0.73%  '
